$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that will no longer be needed (rows 7 through 37),
# collapsing the remaining rows below row 6 (none exist past 37 here).
$ws.Range("A7:A37").EntireRow.Delete()

# Update the consolidated card text for rows 2-6.
$ws.Range("A2").Value = "('Gravecrawler', ['{B}', 'Creature — Zombie', 'Gravecrawler can" + [char]8217 + "t block.', 'You may cast Gravecrawler from your graveyard as long as you control a Zombie.', '2/1'])"
$ws.Range("A3").Value = "('Mondronen Shaman', ['{3}{R}', 'Creature — Human Shaman Werewolf', 'At the beginning of each upkeep, if no spells were cast last turn, transform Mondronen Shaman.', '3/2', " + [char]34 + "Tovolar's Magehunter" + [char]34 + ", 'Creature — Werewolf', 'Whenever an opponent casts a spell, Tovolar" + [char]8217 + "s Magehunter deals 2 damage to that player.', 'At the beginning of each upkeep, if a player cast two or more spells last turn, transform Tovolar" + [char]8217 + "s Magehunter.', '5/5'])"
$ws.Range("A4").Value = "('Ravenous Demon', ['{3}{B}{B}', 'Creature — Demon', 'Sacrifice a Human: Transform Ravenous Demon. Activate this ability only any time you could cast a sorcery.', '4/4', 'Archdemon of Greed', 'Creature — Demon', 'Flying, trample', 'At the beginning of your upkeep, sacrifice a Human. If you can" + [char]8217 + "t, tap Archdemon of Greed and it deals 9 damage to you.', '9/9'])"
$ws.Range("A5").Value = "('Strangleroot Geist', ['{G}{G}', 'Creature — Spirit', 'Haste', 'Undying (When this creature dies, if it had no +1/+1 counters on it, return it to the battlefield under its owner" + [char]8217 + "s control with a +1/+1 counter on it.)', '2/1'])"
$ws.Range("A6").Value = "('Zombie Apocalypse', ['{3}{B}{B}{B}', 'Sorcery', 'Return all Zombie creature cards from your graveyard to the battlefield tapped, then destroy all Humans.'])"
